$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 373.42105
$ws.Range("J33").Value = 624.5
$ws.Range("L33").Value = 624.5
$ws.Range("N33").Value = -1082.5
$ws.Range("H62").Value = 23837490
$ws.Range("I62").Value = 28602988
$ws.Range("K62").Value = 28602988
$ws.Range("M62").Value = -28602364
$ws.Range("H65").Value = 23837490
$ws.Range("I65").Value = 28602988
$ws.Range("K65").Value = 143014940
$ws.Range("M65").Value = -143011820
$ws.Range("H80").Value = 658.3333
$ws.Range("I80").Value = 304.8
$ws.Range("J80").Value = 835.1
$ws.Range("K80").Value = 914.4000000000001
$ws.Range("L80").Value = 2505.3
$ws.Range("M80").Value = 83.59999999999991
$ws.Range("N80").Value = -4501.3
$ws.Range("H83").Value = 658.3333
$ws.Range("I83").Value = 304.8
$ws.Range("J83").Value = 835.1
$ws.Range("K83").Value = 2743.2
$ws.Range("L83").Value = 7515.900000000001
$ws.Range("M83").Value = 2248.8
$ws.Range("N83").Value = -17499.9
$ws.Range("H86").Value = 100001860
$ws.Range("I86").Value = 166668030
$ws.Range("J86").Value = 2588
$ws.Range("K86").Value = 166668030
$ws.Range("L86").Value = 2588
$ws.Range("M86").Value = -166666907
$ws.Range("N86").Value = -4834
$ws.Range("H89").Value = 100001860
$ws.Range("I89").Value = 166668030
$ws.Range("J89").Value = 2588
$ws.Range("K89").Value = 833340150
$ws.Range("L89").Value = 12940
$ws.Range("M89").Value = -833334534
$ws.Range("N89").Value = -24172
$ws.Range("H137").Value = 12058.823
$ws.Range("I137").Value = 8714.700000000001
$ws.Range("K137").Value = 26144.1
$ws.Range("M137").Value = -23594.1
$ws.Range("H138").Value = 2966.9019
$ws.Range("J138").Value = 6886.778
$ws.Range("L138").Value = 20660.334
$ws.Range("N138").Value = -30940.334

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4117.5454
$ws.Range("I2").Value = 2177.7778
$ws.Range("K2").Value = 2177.7778
$ws.Range("M2").Value = -2064.7778
$ws.Range("H97").Value = 1369.55
$ws.Range("I97").Value = 1091.5385
$ws.Range("K97").Value = 1091.5385
$ws.Range("M97").Value = -595.5385000000001
$ws.Range("H102").Value = 5931.727
$ws.Range("I102").Value = 5525
$ws.Range("K102").Value = 5525
$ws.Range("M102").Value = -3903
$ws.Range("H116").Value = 4117.5454
$ws.Range("I116").Value = 2177.7778
$ws.Range("K116").Value = 2177.7778
$ws.Range("M116").Value = 116.2222000000002

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4117.5454
$ws.Range("I3").Value = 2177.7778
$ws.Range("K3").Value = 2177.7778
$ws.Range("M3").Value = -2063.7778
$ws.Range("H94").Value = 47946.76
$ws.Range("I94").Value = 4979.9375
$ws.Range("J94").Value = 124332.22
$ws.Range("K94").Value = 4979.9375
$ws.Range("L94").Value = 124332.22
$ws.Range("M94").Value = -4528.9375
$ws.Range("N94").Value = -125234.22
$ws.Range("H134").Value = 1112553.9
$ws.Range("I134").Value = 1257373
$ws.Range("K134").Value = 3772119
$ws.Range("M134").Value = -3769584

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5291.7915
$ws.Range("I31").Value = 2536.25
$ws.Range("K31").Value = 2536.25
$ws.Range("M31").Value = -2241.25
$ws.Range("H34").Value = 5291.7915
$ws.Range("I34").Value = 2536.25
$ws.Range("K34").Value = 2536.25
$ws.Range("M34").Value = -2334.25
$ws.Range("H99").Value = 6947182
$ws.Range("J99").Value = 3485.625
$ws.Range("L99").Value = 3485.625
$ws.Range("N99").Value = -6481.625
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H126").Value = 6947182
$ws.Range("J126").Value = 3485.625
$ws.Range("L126").Value = 10456.875
$ws.Range("N126").Value = -15396.875

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2375.25
$ws.Range("J22").Value = 2375.25
$ws.Range("L22").Value = 7125.75
$ws.Range("N22").Value = -7463.75
$ws.Range("H23").Value = 848
$ws.Range("I23").Value = 229.66667
$ws.Range("J23").Value = 1466.3334
$ws.Range("K23").Value = 689.00001
$ws.Range("L23").Value = 4399.0002
$ws.Range("M23").Value = -454.00001
$ws.Range("N23").Value = -4869.0002
$ws.Range("H27").Value = 2375.25
$ws.Range("J27").Value = 2375.25
$ws.Range("L27").Value = 7125.75
$ws.Range("N27").Value = -7329.75
$ws.Range("H68").Value = 1251124.8
$ws.Range("J68").Value = 5000000
$ws.Range("L68").Value = 15000000
$ws.Range("N68").Value = -15001622
$ws.Range("H71").Value = 1251124.8
$ws.Range("J71").Value = 5000000
$ws.Range("L71").Value = 45000000
$ws.Range("N71").Value = -45008112

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 16673500
$ws.Range("J29").Value = 500
$ws.Range("L29").Value = 500
$ws.Range("N29").Value = -1080
$ws.Range("H63").Value = 34999.668
$ws.Range("J63").Value = 44999
$ws.Range("L63").Value = 44999
$ws.Range("N63").Value = -46371
$ws.Range("H66").Value = 34999.668
$ws.Range("J66").Value = 44999
$ws.Range("L66").Value = 134997
$ws.Range("N66").Value = -141861
$ws.Range("H122").Value = 4276.324
$ws.Range("I122").Value = 3495.4167
$ws.Range("K122").Value = 10486.2501
$ws.Range("M122").Value = -8036.250100000001

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 75006610
$ws.Range("I136").Value = 31256390
$ws.Range("K136").Value = 93769170
$ws.Range("M136").Value = -93766620

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 9292.714
$ws.Range("I4").Value = 7687.5
$ws.Range("K4").Value = 7687.5
$ws.Range("M4").Value = -7574.5
$ws.Range("H44").Value = 25120
$ws.Range("J44").Value = 25120
$ws.Range("L44").Value = 25120
$ws.Range("N44").Value = -26228
$ws.Range("H107").Value = 527.7273
$ws.Range("I107").Value = 562.2
$ws.Range("K107").Value = 1686.6
$ws.Range("M107").Value = 233.3999999999999
$ws.Range("H132").Value = 9984.875
$ws.Range("I132").Value = 6649.1665
$ws.Range("K132").Value = 19947.4995
$ws.Range("M132").Value = -17417.4995
$ws.Range("H136").Value = 12199940
$ws.Range("I136").Value = 14707094
$ws.Range("K136").Value = 44121282
$ws.Range("M136").Value = -44118732
